$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Experiments")
$ws.Activate()

$ws.Range("A16").Value = "Hu, Ozay et. al. (2018) Revisiting Single Image Depth Estimation"
$ws.Range("B16").Value = 0.866
$ws.Range("C16").Value = 0.975
$ws.Range("D16").Value = 0.993
$ws.Range("F16").Value = 0.53
$ws.Range("G16").Value = 0.115
$ws.Range("I16").Value = 0.05

$ws.Range("B16:D16").NumberFormat = "0.000"
$ws.Range("F16:G16").NumberFormat = "0.000"
$ws.Range("I16").NumberFormat = "0.000"

$ws.Range("J16").Select()
